$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I11").Value = "Each service is a separate app."
$ws.Range("I12").Value = "the auth-service usesthe synchronize mechanism (http) to communicate with the user-service"

$ws.Range("I13").Select()
